$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (price / volume) to stay as plain text so
# numeric-looking strings (e.g. "300.71") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "43.067.16"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "2.304.17"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "300.71"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "98.03"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").Value = "0.521"
$ws.Range("E7").Value = "  +4.30%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("D10").Value = "35.63"
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "17.98"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").Value = "2.663.01"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").Value = "2.362.72"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("D17").Value = "0.786"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "42.975.97"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").Value = "13.38"
$ws.Range("E19").Value = "  +7.93%  "
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").Value = "6.12"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("D22").Value = "68.21"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "238.59"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("D24").Value = "2.18"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").Value = "2.43"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("D28").Value = "168.50"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  -6.66%  "
$ws.Range("D30").Value = "9.13"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "32.97"
$ws.Range("E31").Value = "  -3.72%  "
$ws.Range("D32").Value = "5.18"
$ws.Range("E32").Value = "  +4.46%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "4.81"
$ws.Range("E34").Value = "  +5.32%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "18.10"
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("D41").Value = "2.75"
$ws.Range("E41").Value = "  -2.64%  "
$ws.Range("D42").Value = "2.010.19"
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").Value = "2.16"
$ws.Range("E44").Value = "  -2.70%  "
$ws.Range("D45").Value = "10.17"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "17.42"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("D48").Value = "54.45"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").Value = "2.529.25"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "73.48"
$ws.Range("E51").Value = "  +5.22%  "

# Restore the default (unstyled) cell style now that the values are locked in as text.
$ws.Range("D2:E51").Style = "Normal"

